$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.950.63"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "3.152.72"
$ws.Range("E3").Value = "  +3.09%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.154.25"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("E10").Value = "  +4.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.497"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("E13").Value = "  +13.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("D15").Value = "3.671.44"
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "65.077.93"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "3.156.32"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "506.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.08%  "
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.63%  "
$ws.Range("E30").Value = "  +7.31%  "
$ws.Range("E31").Value = "  +4.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0898"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "463.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.61%  "
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("D42").Value = "3.050.75"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.42%  "
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("D47").Value = "0.0₃0581"
$ws.Range("E47").Value = "  +12.18%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.70%  "
